# Applies the "improved random forest model (dog-to-human)" edit:
#  - file_description sheet: add two new rows describing the transposed_OTU
#    and metadata_otu_merged_non-rarefied files
#  - metadata column dictionary sheet: delete the empty "Description" column
#  - window/view cosmetics

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: file_description ----
$ws1 = $wb.Worksheets.Item("file_description")
$ws1.Range("B7").Value = "biom_analysis"
$ws1.Range("B8").Value = "biom_analysis"
$ws1.Range("C7").Value = "transposed_OTU"
$ws1.Range("C8").Value = "metadata_otu_merged_non-rarefied"
$ws1.Range("D7").Value = "original out_table, transposed (non-rarefied)"
$ws1.Range("D8").Value = "transposed_OTU merged with metadata"
$ws1.Range("D9").Select()

# ---- Sheet 2: metadata column dictionary ----
$ws2 = $wb.Worksheets.Item("metadata column dictionary")
$ws2.Columns.Item(3).Delete()
$ws2.Range("C1:C1048576").Select()
